$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.178.56"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").Value = "'1.591.36"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'211.98"
$ws.Range("E5").Value = "  +0.87%  "
$ws.Range("D6").Value = "'0.502"
$ws.Range("E6").Value = "  -0.91%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Value = "'0.0605"
$ws.Range("E9").Value = "  -0.87%  "
$ws.Range("D10").Value = "'18.97"
$ws.Range("E10").Value = "  -2.30%  "
$ws.Range("D11").Value = "'0.0846"
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("D12").Value = "'1.812.56"
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("D13").Value = "'1.588.38"
$ws.Range("E13").Value = "  -0.02%  "
$ws.Range("E14").Value = "  -1.39%  "
$ws.Range("E15").Value = "  -1.98%  "
$ws.Range("D16").Value = "'63.59"
$ws.Range("E16").Value = "  -1.22%  "
$ws.Range("D17").Value = "'26.172.89"
$ws.Range("E17").Value = "  -0.56%  "
$ws.Range("D18").Value = "'0.0₃0725"
$ws.Range("E18").Value = "  -0.48%  "
$ws.Range("D19").Value = "'213.99"
$ws.Range("E19").Value = "  +1.46%  "
$ws.Range("E20").Value = "  -1.57%  "
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("D22").Value = "'4.25"
$ws.Range("E22").Value = "  -0.59%  "
$ws.Range("D23").Value = "'9.03"
$ws.Range("E23").Value = "  +0.79%  "
$ws.Range("D24").Value = "'2.11"
$ws.Range("E24").Value = "  -1.79%  "
$ws.Range("D25").Value = "'144.88"
$ws.Range("E25").Value = "  +0.19%  "
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("D27").Value = "'6.96"
$ws.Range("E27").Value = "  -1.21%  "
$ws.Range("D28").Value = "'0.112"
$ws.Range("E28").Value = "  -1.35%  "
$ws.Range("D29").Value = "'15.07"
$ws.Range("E29").Value = "  -1.09%  "
$ws.Range("E30").Value = "  -2.59%  "
$ws.Range("E31").Value = "  +0.22%  "
$ws.Range("D32").Value = "'3.16"
$ws.Range("E32").Value = "  -1.82%  "
$ws.Range("D33").Value = "'1.419.85"
$ws.Range("E33").Value = "  +7.88%  "
$ws.Range("D34").Value = "'2.95"
$ws.Range("E34").Value = "  -1.63%  "
$ws.Range("D35").Value = "'2.42"
$ws.Range("E35").Value = "  -0.85%  "
$ws.Range("E36").Value = "  -0.83%  "
$ws.Range("E37").Value = "  -3.79%  "
$ws.Range("E38").Value = "  -1.60%  "
$ws.Range("D39").Value = "'5.90"
$ws.Range("E39").Value = "  +4.62%  "
$ws.Range("D40").Value = "'0.823"
$ws.Range("E40").Value = "  +2.18%  "
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("D42").Value = "'0.967"
$ws.Range("E42").Value = "  -10.16%  "
$ws.Range("D43").Value = "'0.765"
$ws.Range("E43").Value = "  -0.24%  "
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("D45").Value = "'1.725.23"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").Value = "'60.95"
$ws.Range("E46").Value = "  -2.14%  "
$ws.Range("D47").Value = "'86.96"
$ws.Range("E47").Value = "  -0.63%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "'0.0₆0103"
$ws.Range("E48").Value = "  -1.09%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'1.48"
$ws.Range("E49").Value = "  -0.28%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.0502"
$ws.Range("E50").Value = "  -0.92%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.0957"
$ws.Range("E51").Value = "  -1.94%  "
